# "MIRADOR OPTICO" price list refresh:
#  - bump the quote date in A1 by one month (45406 -> 45436, i.e. 24/04/2024 -> 24/05/2024)
#  - update the two product prices in D22/D23 (1497.908 -> 2950.798)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# A1 is formatted as a date (numFmt 14); set an exact midnight date so no time
# fraction gets written to the serial value.
$ws.Range("A1").Value = (Get-Date -Year 2024 -Month 5 -Day 24 -Hour 0 -Minute 0 -Second 0)

# New prices for "MIRADOR OPTICO Hº Niquel." (row 22) and "MIRADOR OPTICO Hº Bceado." (row 23)
$ws.Range("D22").Value = 2950.798
$ws.Range("D23").Value = 2950.798

# Re-touch the B:C merges of rows 22/23 (unmerge+remerge, swapping the order in
# which they're re-created) so the workbook's merged-cell bookkeeping ends up
# ordered the same way as in the refreshed file. Formatting is stashed in a
# scratch area first and pasted back afterwards, since unmerging resets it.
$ws.Range("B22:C22").Copy()
$ws.Range("H1:I1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("B23:C23").Copy()
$ws.Range("H2:I2").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("B22:C22").UnMerge()
$ws.Range("B23:C23").UnMerge()
$ws.Range("B23:C23").Merge()
$ws.Range("B22:C22").Merge()

$ws.Range("H1:I1").Copy()
$ws.Range("B22:C22").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("H2:I2").Copy()
$ws.Range("B23:C23").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("H1:I2").Clear()
